# 9th Stab - Cosmetic Changes
# Insert two new weekly date columns ("Jun_17" and "Jun_15") to the left of
# the existing most-recent-week column (old column C, "Jun_10" data), so the
# report now reads, left to right: Ticker | Jun_17 | Jun_15 | Jun_13 | Jun_10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Jun_10" column (C) two columns to the right (-> E),
# opening up fresh columns C and D for the two new weeks.
$ws.Columns("C:D").Insert()

# New header row: newest week first.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Every data row gets the "UN" (unknown/no rating change) placeholder in the
# two freshly inserted week columns, matching the placeholder already used
# throughout column B.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Keep the cosmetic column widths consistent across the now-adjacent date
# columns (matches the pre-existing 8-character-wide column formatting).
$ws.Columns("C").ColumnWidth = 7.166666666666667
$ws.Columns("D").ColumnWidth = 7.166666666666667
$ws.Columns("E").ColumnWidth = 7.166666666666667
